$wb = $excel.ActiveWorkbook

# Select Sheet1 and set its view state (selection) before adding the new sheet
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("C65").Select() | Out-Null

# Add a new worksheet "Sheet2" positioned right after "Sheet1" (becomes the last sheet)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws.Name = "Sheet2"

# Match the page setup of the authored sheet (margins in points: 0.75in/0.75in/1in/1in/0.5in/0.5in)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
$ws.PageSetup.Orientation = 1

# Seed shared-string insertion order to match the authored workbook
$ws.Range("A1").Value = "design"
$ws.Range("E1").Value = "requirement"
$ws.Range("A2").Value = "F1 measure"
$ws.Range("A16").Value = "Precision"
$ws.Range("A29").Value = "Recall"
$ws.Range("A3").Value = "LR"
$ws.Range("B3").Value = "Binary"
$ws.Range("C3").Value = "Bayes"

# Remaining cell values/formulas
$ws.Range("E2").Value = "F1 measure"
$ws.Range("E3").Value = "LR"
$ws.Range("F3").Value = "Binary"
$ws.Range("G3").Value = "Bayes"
$ws.Range("A4").Value = 0.51700000000000002
$ws.Range("B4").Value = 0.56299999999999994
$ws.Range("C4").Value = 0.13400000000000001
$ws.Range("E4").Value = 0.36399999999999999
$ws.Range("F4").Value = 0.41399999999999998
$ws.Range("G4").Value = 0.03
$ws.Range("A5").Value = 0.73099999999999998
$ws.Range("B5").Value = 0.82199999999999995
$ws.Range("C5").Value = 0.52500000000000002
$ws.Range("E5").Value = 0.255
$ws.Range("F5").Value = 0.76
$ws.Range("G5").Value = 0.33500000000000002
$ws.Range("A6").Value = 0.81399999999999995
$ws.Range("B6").Value = 0.627
$ws.Range("C6").Value = 0.29399999999999998
$ws.Range("E6").Value = 0.76
$ws.Range("F6").Value = 0.93400000000000005
$ws.Range("G6").Value = 0.189
$ws.Range("A7").Value = 0.60099999999999998
$ws.Range("B7").Value = 0.48799999999999999
$ws.Range("C7").Value = 0.106
$ws.Range("E7").Value = 0.93400000000000005
$ws.Range("F7").Value = 0.38100000000000001
$ws.Range("G7").Value = 0.02
$ws.Range("A8").Value = 0.47
$ws.Range("B8").Value = 0.76700000000000002
$ws.Range("C8").Value = 0.435
$ws.Range("E8").Value = 0.38100000000000001
$ws.Range("F8").Value = 0.48499999999999999
$ws.Range("G8").Value = 0.097000000000000003
$ws.Range("A9").Value = 0.74399999999999999
$ws.Range("B9").Value = 0.48
$ws.Range("C9").Value = 0.35299999999999998
$ws.Range("E9").Value = 0.47599999999999998
$ws.Range("F9").Value = 0.095000000000000001
$ws.Range("G9").Value = 0.02
$ws.Range("A10").Value = 0.50900000000000001
$ws.Range("B10").Value = 0.495
$ws.Range("C10").Value = 0.35
$ws.Range("E10").Value = 0.090999999999999998
$ws.Range("F10").Value = 0.45900000000000002
$ws.Range("G10").Value = 0.039
$ws.Range("A11").Value = 0.49199999999999999
$ws.Range("B11").Value = 0.73699999999999999
$ws.Range("C11").Value = 0.224
$ws.Range("E11").Value = 0.5
$ws.Range("F11").Value = 0.28599999999999998
$ws.Range("G11").Value = 0.029000000000000001
$ws.Range("A12").Value = 0.78300000000000003
$ws.Range("B12").Value = 0.81100000000000005
$ws.Range("C12").Value = 0.42899999999999999
$ws.Range("E12").Value = 0.46200000000000002
$ws.Range("F12").Value = 0.46800000000000003
$ws.Range("G12").Value = 0.13700000000000001
$ws.Range("A13").Value = 0.54
$ws.Range("B13").Value = 0.55800000000000005
$ws.Range("C13").Value = 0.23300000000000001
$ws.Range("E13").Value = 0.875
$ws.Range("F13").Value = 0.82
$ws.Range("G13").Value = 0.121
$ws.Range("A14").Formula = "=AVERAGE(A4:A13)"
$ws.Range("B14").Formula = "=AVERAGE(B4:B13)"
$ws.Range("C14").Formula = "=AVERAGE(C4:C13)"
$ws.Range("E14").Formula = "=AVERAGE(E4:E13)"
$ws.Range("F14").Formula = "=AVERAGE(F4:F13)"
$ws.Range("G14").Formula = "=AVERAGE(G4:G13)"
$ws.Range("E16").Value = "Precision"
$ws.Range("A17").Value = "LR"
$ws.Range("B17").Value = "Binary"
$ws.Range("C17").Value = "Bayes"
$ws.Range("E17").Value = "LR"
$ws.Range("F17").Value = "Binary"
$ws.Range("G17").Value = "Bayes"
$ws.Range("A18").Value = 0.55400000000000005
$ws.Range("B18").Value = 0.62
$ws.Range("C18").Value = 0.071999999999999995
$ws.Range("E18").Value = 0.35299999999999998
$ws.Range("F18").Value = 0.46200000000000002
$ws.Range("G18").Value = 0.014999999999999999
$ws.Range("A19").Value = 0.80800000000000005
$ws.Range("B19").Value = 0.79
$ws.Range("C19").Value = 0.35799999999999998
$ws.Range("E19").Value = 0.16700000000000001
$ws.Range("F19").Value = 0.79200000000000004
$ws.Range("G19").Value = 0.20699999999999999
$ws.Range("A20").Value = 0.78800000000000003
$ws.Range("B20").Value = 0.84
$ws.Range("C20").Value = 0.18099999999999999
$ws.Range("E20").Value = 0.80200000000000005
$ws.Range("F20").Value = 0.91400000000000003
$ws.Range("G20").Value = 0.105
$ws.Range("A21").Value = 0.79200000000000004
$ws.Range("B21").Value = 0.63300000000000001
$ws.Range("C21").Value = 0.057000000000000002
$ws.Range("E21").Value = 0.91400000000000003
$ws.Range("F21").Value = 0.8
$ws.Range("G21").Value = 0.01
$ws.Range("A22").Value = 0.57399999999999995
$ws.Range("B22").Value = 0.89500000000000002
$ws.Range("C22").Value = 0.28799999999999998
$ws.Range("E22").Value = 0.8
$ws.Range("F22").Value = 0.64100000000000001
$ws.Range("G22").Value = 0.050999999999999997
$ws.Range("A23").Value = 0.877
$ws.Range("B23").Value = 0.80700000000000005
$ws.Range("C23").Value = 0.22700000000000001
$ws.Range("E23").Value = 0.61
$ws.Range("F23").Value = 0.14299999999999999
$ws.Range("G23").Value = 0.01
$ws.Range("A24").Value = 0.77900000000000003
$ws.Range("B24").Value = 0.65800000000000003
$ws.Range("C24").Value = 0.224
$ws.Range("E24").Value = 0.125
$ws.Range("F24").Value = 0.34699999999999998
$ws.Range("G24").Value = 0.02
$ws.Range("A25").Value = 0.64600000000000002
$ws.Range("B25").Value = 0.81899999999999995
$ws.Range("C25").Value = 0.14000000000000001
$ws.Range("E25").Value = 0.373
$ws.Range("F25").Value = 0.19400000000000001
$ws.Range("G25").Value = 0.014999999999999999
$ws.Range("A26").Value = 0.79800000000000004
$ws.Range("B26").Value = 0.81499999999999995
$ws.Range("C26").Value = 0.27500000000000002
$ws.Range("E26").Value = 0.70899999999999996
$ws.Range("F26").Value = 0.70199999999999996
$ws.Range("G26").Value = 0.073999999999999996
$ws.Range("A27").Value = 0.54400000000000004
$ws.Range("B27").Value = 0.56699999999999995
$ws.Range("C27").Value = 0.13300000000000001
$ws.Range("E27").Value = 0.91300000000000003
$ws.Range("F27").Value = 0.80600000000000005
$ws.Range("G27").Value = 0.064000000000000001
$ws.Range("A28").Formula = "=AVERAGE(A18:A27)"
$ws.Range("B28").Formula = "=AVERAGE(B18:B27)"
$ws.Range("C28").Formula = "=AVERAGE(C18:C27)"
$ws.Range("E28").Formula = "=AVERAGE(E18:E27)"
$ws.Range("F28").Formula = "=AVERAGE(F18:F27)"
$ws.Range("G28").Formula = "=AVERAGE(G18:G27)"
$ws.Range("E29").Value = "Recall"
$ws.Range("A30").Value = "LR"
$ws.Range("B30").Value = "Binary"
$ws.Range("C30").Value = "Bayes"
$ws.Range("E30").Value = "LR"
$ws.Range("F30").Value = "Binary"
$ws.Range("G30").Value = "Bayes"
$ws.Range("A31").Value = 0.48399999999999999
$ws.Range("B31").Value = 0.51600000000000001
$ws.Range("C31").Value = 0.874
$ws.Range("E31").Value = 0.375
$ws.Range("F31").Value = 0.375
$ws.Range("G31").Value = 0.86799999999999999
$ws.Range("A32").Value = 0.66800000000000004
$ws.Range("B32").Value = 0.85799999999999998
$ws.Range("C32").Value = 0.98499999999999999
$ws.Range("E32").Value = 0.54500000000000004
$ws.Range("F32").Value = 0.73099999999999998
$ws.Range("G32").Value = 0.877
$ws.Range("A33").Value = 0.84299999999999997
$ws.Range("B33").Value = 0.5
$ws.Range("C33").Value = 0.78600000000000003
$ws.Range("E33").Value = 0.72199999999999998
$ws.Range("F33").Value = 0.95499999999999996
$ws.Range("G33").Value = 0.97799999999999998
$ws.Range("A34").Value = 0.48399999999999999
$ws.Range("B34").Value = 0.39700000000000002
$ws.Range("C34").Value = 0.872
$ws.Range("E34").Value = 0.95499999999999996
$ws.Range("F34").Value = 0.25
$ws.Range("G34").Value = 1
$ws.Range("A35").Value = 0.39700000000000002
$ws.Range("B35").Value = 0.67
$ws.Range("C35").Value = 0.89
$ws.Range("E35").Value = 0.25
$ws.Range("F35").Value = 0.39100000000000001
$ws.Range("G35").Value = 0.84399999999999997
$ws.Range("A36").Value = 0.64500000000000002
$ws.Range("B36").Value = 0.34200000000000003
$ws.Range("C36").Value = 0.79100000000000004
$ws.Range("E36").Value = 0.39100000000000001
$ws.Range("F36").Value = 0.070999999999999994
$ws.Range("G36").Value = 0.78600000000000003
$ws.Range("A37").Value = 0.378
$ws.Range("B37").Value = 0.39700000000000002
$ws.Range("C37").Value = 0.80100000000000005
$ws.Range("E37").Value = 0.070999999999999994
$ws.Range("F37").Value = 0.68
$ws.Range("G37").Value = 0.84
$ws.Range("A38").Value = 0.39700000000000002
$ws.Range("B38").Value = 0.67100000000000004
$ws.Range("C38").Value = 0.56000000000000005
$ws.Range("E38").Value = 0.76
$ws.Range("F38").Value = 0.54500000000000004
$ws.Range("G38").Value = 0.95499999999999996
$ws.Range("A39").Value = 0.77
$ws.Range("B39").Value = 0.80800000000000005
$ws.Range("C39").Value = 0.97099999999999997
$ws.Range("E39").Value = 0.34200000000000003
$ws.Range("F39").Value = 0.35099999999999998
$ws.Range("G39").Value = 0.90400000000000003
$ws.Range("A40").Value = 0.53600000000000003
$ws.Range("B40").Value = 0.55000000000000004
$ws.Range("C40").Value = 0.94699999999999995
$ws.Range("E40").Value = 0.84
$ws.Range("F40").Value = 0.83299999999999996
$ws.Range("G40").Value = 0.96699999999999997
$ws.Range("A41").Formula = "=AVERAGE(A31:A40)"
$ws.Range("B41").Formula = "=AVERAGE(B31:B40)"
$ws.Range("C41").Formula = "=AVERAGE(C31:C40)"
$ws.Range("E41").Formula = "=AVERAGE(E31:E40)"
$ws.Range("F41").Formula = "=AVERAGE(F31:F40)"
$ws.Range("G41").Formula = "=AVERAGE(G31:G40)"

$ws.Range("H43").Select() | Out-Null
